# changed index order provincie
# Swap the "Provincienaam" (province name) and "Provinciecode" (province code)
# columns: column A now holds the code, column B now holds the name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header labels in row 1 (A1 <-> B1)
$headerA = $ws.Cells.Item(1, 1).Value()
$headerB = $ws.Cells.Item(1, 2).Value()
$ws.Cells.Item(1, 1).Value = $headerB
$ws.Cells.Item(1, 2).Value = $headerA

# Swap province name (A) and province code (B) for every data row (2-13)
for ($r = 2; $r -le 13; $r++) {
    $nameValue = $ws.Cells.Item($r, 1).Value()
    $codeValue = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($r, 1).Value = $codeValue
    $ws.Cells.Item($r, 2).Value = $nameValue
}
